# Generate Report for Archive
# - Status text moves from "Ready for handoff" to "In Translation" across
#   the Overview / zh-cn / de-de sheets.
# - The Overview "zh-cn"/"de-de" status columns (E, F) and the per-locale
#   "Status" column (C) are narrower now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Replace($oldStatus, $newStatus) | Out-Null
}

# Narrow the status columns to match the shorter text. Columns("...").ColumnWidth
# snaps to Excel's internal pixel grid for the Normal-style font, so we pick
# the closest reachable width to the refreshed auto-fit value (~13.41 chars).
$newWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns("E").ColumnWidth = $newWidth
$wsOverview.Columns("F").ColumnWidth = $newWidth

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns("C").ColumnWidth = $newWidth

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns("C").ColumnWidth = $newWidth
